# The commit clears the "assignment name/score" cells in column F for rows
# 6-9 (an assignment/homework entry that was deleted). In the target OOXML
# these cells stay present as string cells that point at the (already
# existing) empty shared string, rather than being removed outright.
#
# A plain `Value = ""` (or ClearContents) deletes the cell entirely in this
# engine, so instead we write a value that is guaranteed to resolve to an
# empty string (a lone `'` is Excel's "treat next text as text" prefix and
# is not part of the stored value) and then strip the incidental
# quote-prefix formatting that operation introduces, leaving a bare empty
# string cell behind - matching F6/F7/F8/F9 -> "" in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("F6", "F7", "F8", "F9")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value = "'"
    $rng.ClearFormats()
}
